# Actualizacion Datos Personales 4 nov
#
# Updates the "Rescatables" sheet: two new students (HUESCA/GARCIA/ALDAIR OMAR
# and VAZQUEZ/VICTORIANO/MARIAN) are inserted into the roster, and the table
# grows from 17 to 19 data rows (A1:G18 -> A1:G20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Final roster, in row order (row 2 .. row 20).
#        Row  NC(A)            Paterno(B)    Materno(C)    Nombres(D)           Grupo(F)
$data = @(
  @(2,  20330051920002, "ADELL",      "AGUILAR",     "CRISTOPHER ALAIN", "3AEM"),
  @(3,  20330051920361, "BARRAGAN",   "CASTILLO",    "JAIR",             "3AEM"),
  @(4,  20330051920245, "PONCE",      "GOMEZ",       "ALETHIA LUCIA",    "3ALCM"),
  @(5,  20330051920168, "HUESCA",     "GARCIA",      "ALDAIR OMAR",      "3ARHM"),
  @(6,  20330051920184, "VAZQUEZ",    "VICTORIANO",  "MARIAN",           "3ARHM"),
  @(7,  20330051920042, "ANASTACIO",  "ROMERO",      "HIRAM FABIAN",     "3BEM"),
  @(8,  20330051920011, "ESPIRITU",   "TEQUIHUATLE", "ALEJANDRO",        "3AEM"),
  @(9,  20330051920012, "ESPIRITU",   "TEQUIHUATLE", "IGNACIO",          "3AEM"),
  @(10, 20330051920027, "ROMAN",      "GONZALEZ",    "LUIS ALEJANDRO",   "3AEM"),
  @(11, 20330051920036, "XOTLANIHUA", "RODRIGUEZ",   "JOSE ANTONIO",     "3AEM"),
  @(12, 20330051920224, "COLOHUA",    "RAMIREZ",     "FERNANDA",         "3ALCM"),
  @(13, 20330051920228, "HERNANDEZ",  "GALEOTE",     "GERMAN ISAI",      "3ALCM"),
  @(14, 20330051920158, "BAUTISTA",   "DIAZ",        "DINA BERENICE",    "3ARHM"),
  @(15, 20330051920161, "CANSECO",    "LEAL",        "ANGELA",           "3ARHM"),
  @(16, 20330051920172, "MATA",       "CANSECO",     "CRISTIAN ARTURO",  "3ARHM"),
  @(17, 20330051920178, "QUIRIZ",     "RAMOS",       "MONICA",           "3ARHM"),
  @(18, 20330051920254, "SANCHEZ",    "ROMERO",      "BERENICE",         "3ARHM"),
  @(19, 20330051920046, "CID",        "VALENCIA",    "JESUS",            "3BEM"),
  @(20, 20330051920111, "VALENTE",    "GAMEZ",       "ABIUD",            "3BEM")
)

# NC (matricula numbers), column A
foreach ($rec in $data) {
    $ws.Range("A" + $rec[0]).Value = $rec[1]
}

# Paterno, column B
foreach ($rec in $data) {
    $ws.Range("B" + $rec[0]).Value = $rec[2]
}

# Materno, column C
foreach ($rec in $data) {
    $ws.Range("C" + $rec[0]).Value = $rec[3]
}

# Nombres, column D
foreach ($rec in $data) {
    $ws.Range("D" + $rec[0]).Value = $rec[4]
}

# Nombre_Largo (subject), column E - same for every row
foreach ($rec in $data) {
    $ws.Range("E" + $rec[0]).Value = "GEOMETRÍA ANALÍTICA"
}

# Grupo, column F
foreach ($rec in $data) {
    $ws.Range("F" + $rec[0]).Value = $rec[5]
}

# Reprobadas, column G - constant 6 for every row
foreach ($rec in $data) {
    $ws.Range("G" + $rec[0]).Value = 6
}
